# Rename worksheets to their lowercase (snake_case) equivalents,
# dropping the leading "d_"/"f_" domain/fact prefix.

$wb = $excel.ActiveWorkbook

$renameMap = @{
    "d_ANIO_REGISTRO_BDD"         = "anio_registro_BDD"
    "d_MES_REGISTRO_BDD"          = "mes_registro_BDD"
    "d_MES_INCIDENTE_BDD"         = "mes_incidente_BDD"
    "d_ARTICULO_BDD"              = "articulo_BDD"
    "d_DELITO_BDD"                = "delito_BDD"
    "d_DELITO_CIRCUNSTANCIAL_BDD" = "delito_circunstancial_BDD"
    "d_ESTADO_PROCESAL_BDD"       = "estado_procesal_BDD"
    "d_ETAPA_ACTUAL_BDD"          = "etapa_actual_BDD"
    "d_PROVINCIA_INCIDENTE_BDD"   = "provincia_incidente_BDD"
    "d_CANTON_INCIDENTE_BDD"      = "canton_incidente_BDD"
    "d_TIPO_DELITO_BDD"           = "tipo_delito_BDD"
    "d_TIPO_FLAGRANTE_BDD"        = "tipo_flagrante_BDD"
    "d_ESTADO_NDD_BDD"            = "estado_ndd_BDD"
    "f_GRUPO_HORAINC_BDD"         = "grupo_horainc_BDD"
}

foreach ($ws in $wb.Worksheets) {
    $oldName = $ws.Name
    if ($renameMap.ContainsKey($oldName)) {
        $ws.Name = $renameMap[$oldName]
    }
}
